$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 428.73334
$ws.Range("I19").Value = 263.75
$ws.Range("J19").Value = 488.72726
$ws.Range("K19").Value = 263.75
$ws.Range("L19").Value = 488.72726
$ws.Range("M19").Value = -88.75
$ws.Range("N19").Value = -838.72726
$ws.Range("H98").Value = 2057.2727
$ws.Range("I98").Value = 2081.111
$ws.Range("J98").Value = 1950
$ws.Range("K98").Value = 2081.111
$ws.Range("L98").Value = 1950
$ws.Range("M98").Value = -583.1109999999999
$ws.Range("N98").Value = -4946
$ws.Range("H121").Value = 2433.625
$ws.Range("I121").Value = 890
$ws.Range("J121").Value = 2654.1428
$ws.Range("K121").Value = 2670
$ws.Range("L121").Value = 7962.428400000001
$ws.Range("M121").Value = -923
$ws.Range("N121").Value = -11456.4284
$ws.Range("H122").Value = 2057.2727
$ws.Range("I122").Value = 2081.111
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 6243.333
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -3793.333
$ws.Range("N122").Value = -10750
$ws.Range("H131").Value = 2299.7
$ws.Range("I131").Value = 1123.75
$ws.Range("J131").Value = 2593.6875
$ws.Range("K131").Value = 3371.25
$ws.Range("L131").Value = 7781.0625
$ws.Range("M131").Value = 1668.75
$ws.Range("N131").Value = -17861.0625
$ws.Range("H141").Value = 2156.923
$ws.Range("I141").Value = 1703.8096
$ws.Range("J141").Value = 4060
$ws.Range("K141").Value = 5111.4288
$ws.Range("L141").Value = 12180
$ws.Range("M141").Value = 68.57120000000032
$ws.Range("N141").Value = -22540

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6459.697
$ws.Range("I32").Value = 4330.6123
$ws.Range("K32").Value = 4330.6123
$ws.Range("M32").Value = -4043.6123
$ws.Range("H110").Value = 1580.6842
$ws.Range("I110").Value = 1403.6666
$ws.Range("J110").Value = 1884.1428
$ws.Range("K110").Value = 1403.6666
$ws.Range("L110").Value = 1884.1428
$ws.Range("M110").Value = 641.3334
$ws.Range("N110").Value = -5974.1428
$ws.Range("H132").Value = 2652.8223
$ws.Range("I132").Value = 1618.1818
$ws.Range("J132").Value = 5498.0835
$ws.Range("K132").Value = 4854.5454
$ws.Range("L132").Value = 16494.2505
$ws.Range("M132").Value = -2324.5454
$ws.Range("N132").Value = -21554.2505

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H94").Value = 994.1905
$ws.Range("I94").Value = 472.6
$ws.Range("J94").Value = 2298.1667
$ws.Range("K94").Value = 472.6
$ws.Range("L94").Value = 2298.1667
$ws.Range("M94").Value = -21.60000000000002
$ws.Range("N94").Value = -3200.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3100
$ws.Range("I86").Value = 3300
$ws.Range("K86").Value = 3300
$ws.Range("M86").Value = -2177
$ws.Range("H89").Value = 3100
$ws.Range("I89").Value = 3300
$ws.Range("K89").Value = 16500
$ws.Range("M89").Value = -10884
$ws.Range("H94").Value = 3681.658
$ws.Range("I94").Value = 3780.1428
$ws.Range("J94").Value = 3624.2083
$ws.Range("K94").Value = 3780.1428
$ws.Range("L94").Value = 3624.2083
$ws.Range("M94").Value = -3329.1428
$ws.Range("N94").Value = -4526.2083
$ws.Range("H99").Value = 6585360.5
$ws.Range("I99").Value = 12065.777
$ws.Range("J99").Value = 12501326
$ws.Range("K99").Value = 12065.777
$ws.Range("L99").Value = 12501326
$ws.Range("M99").Value = -10567.777
$ws.Range("N99").Value = -12504322
$ws.Range("H126").Value = 6585360.5
$ws.Range("I126").Value = 12065.777
$ws.Range("J126").Value = 12501326
$ws.Range("K126").Value = 36197.331
$ws.Range("L126").Value = 37503978
$ws.Range("M126").Value = -33727.331
$ws.Range("N126").Value = -37508918
$ws.Range("H134").Value = 3022.05
$ws.Range("I134").Value = 2956.257
$ws.Range("K134").Value = 8868.771000000001
$ws.Range("M134").Value = -6333.771000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 3600
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3600
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 10800
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -11656
$ws.Range("H91").Value = 3600
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3600
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 10800
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -13764
$ws.Range("H92").Value = 614.8333
$ws.Range("I92").Value = 400
$ws.Range("J92").Value = 657.8
$ws.Range("K92").Value = 1200
$ws.Range("L92").Value = 1973.4
$ws.Range("M92").Value = 48
$ws.Range("N92").Value = -4469.4
$ws.Range("H113").Value = 984095.5
$ws.Range("I113").Value = 1087446.1
$ws.Range("J113").Value = 667153.6
$ws.Range("K113").Value = 3262338.3
$ws.Range("L113").Value = 2001460.8
$ws.Range("M113").Value = -3260168.3
$ws.Range("N113").Value = -2005800.8
$ws.Range("H131").Value = 1852835.2
$ws.Range("I131").Value = 12500454
$ws.Range("J131").Value = 1075.4783
$ws.Range("K131").Value = 37501362
$ws.Range("L131").Value = 3226.4349
$ws.Range("M131").Value = -37496322
$ws.Range("N131").Value = -13306.4349
$ws.Range("H134").Value = 8950.075000000001
$ws.Range("I134").Value = 7645.294
$ws.Range("J134").Value = 9914.478999999999
$ws.Range("K134").Value = 22935.882
$ws.Range("L134").Value = 29743.437
$ws.Range("M134").Value = -17865.882
$ws.Range("N134").Value = -39883.437

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1820.8948
$ws.Range("I113").Value = 1209.7
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1209.7
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 960.3
$ws.Range("N113").Value = -6840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 7210
$ws.Range("J24").Value = 7210
$ws.Range("L24").Value = 7210
$ws.Range("N24").Value = -7670
$ws.Range("H96").Value = 926.9355
$ws.Range("I96").Value = 949.6667
$ws.Range("K96").Value = 949.6667
$ws.Range("M96").Value = 423.3333
